$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D (Price), E (Volume 1h), and G (Hora) keep their original
# text storage (these columns hold numeric-looking / percentage-looking
# strings that must remain literal text, not be reinterpreted as numbers).
$ws.Range("D2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "303.26"
$ws.Range("E2").Value = "5.31%"
$ws.Range("G2").Value = "12"
$ws.Range("D3").Value = "35.37"
$ws.Range("E3").Value = "14.07%"
$ws.Range("G3").Value = "12"
$ws.Range("D4").Value = "5.185"
$ws.Range("E4").Value = "4.54%"
$ws.Range("G4").Value = "12"
$ws.Range("D5").Value = "0.07829"
$ws.Range("E5").Value = "6.51%"
$ws.Range("G5").Value = "12"
$ws.Range("D6").Value = "2.319"
$ws.Range("E6").Value = "0.17%"
$ws.Range("G6").Value = "12"
$ws.Range("D7").Value = "7.991"
$ws.Range("E7").Value = "3.50%"
$ws.Range("G7").Value = "12"
$ws.Range("D8").Value = "3.976"
$ws.Range("E8").Value = "6.87%"
$ws.Range("G8").Value = "12"
$ws.Range("D9").Value = "0.9299"
$ws.Range("E9").Value = "2.74%"
$ws.Range("G9").Value = "12"
$ws.Range("D10").Value = "0.1005"
$ws.Range("E10").Value = "10.19%"
$ws.Range("G10").Value = "12"
$ws.Range("D11").Value = "0.1836"
$ws.Range("E11").Value = "8.70%"
$ws.Range("G11").Value = "12"
$ws.Range("D12").Value = "0.08537"
$ws.Range("E12").Value = "3.52%"
$ws.Range("G12").Value = "12"
$ws.Range("D13").Value = "0.03396"
$ws.Range("E13").Value = "8.46%"
$ws.Range("G13").Value = "12"
$ws.Range("D14").Value = "0.09897"
$ws.Range("G14").Value = "12"
$ws.Range("D15").Value = "0.001483"
$ws.Range("E15").Value = "-0.89%"
$ws.Range("G15").Value = "12"
$ws.Range("D16").Value = "0.04636"
$ws.Range("E16").Value = "2.52%"
$ws.Range("G16").Value = "12"
$ws.Range("D17").Value = "0.005725"
$ws.Range("E17").Value = "-1.28%"
$ws.Range("G17").Value = "12"
$ws.Range("D18").Value = "3.472"
$ws.Range("E18").Value = "-0.71%"
$ws.Range("G18").Value = "12"
$ws.Range("D19").Value = "2.102"
$ws.Range("E19").Value = "0.26%"
$ws.Range("G19").Value = "12"
$ws.Range("D20").Value = "0.3417"
$ws.Range("E20").Value = "2.84%"
$ws.Range("G20").Value = "12"
$ws.Range("D21").Value = "0.1323"
$ws.Range("E21").Value = "1.82%"
$ws.Range("G21").Value = "12"
$ws.Range("D22").Value = "4.557"
$ws.Range("E22").Value = "8.65%"
$ws.Range("G22").Value = "12"
$ws.Range("D23").Value = "0.2383"
$ws.Range("E23").Value = "11.91%"
$ws.Range("G23").Value = "12"
$ws.Range("D24").Value = "0.001219"
$ws.Range("E24").Value = "0.98%"
$ws.Range("G24").Value = "12"
$ws.Range("D25").Value = "0.004441"
$ws.Range("E25").Value = "6.58%"
$ws.Range("G25").Value = "12"
$ws.Range("D26").Value = "0.0001298"
$ws.Range("E26").Value = "-0.20%"
$ws.Range("G26").Value = "12"
$ws.Range("D27").Value = "0.0003391"
$ws.Range("E27").Value = "0.01%"
$ws.Range("G27").Value = "12"
$ws.Range("G28").Value = "12"
$ws.Range("G29").Value = "12"
$ws.Range("G30").Value = "12"
$ws.Range("G31").Value = "12"
$ws.Range("G32").Value = "12"
$ws.Range("G33").Value = "12"
$ws.Range("G34").Value = "12"
$ws.Range("G35").Value = "12"
$ws.Range("G36").Value = "12"
$ws.Range("G37").Value = "12"
$ws.Range("G38").Value = "12"
$ws.Range("D39").Value = "0.01758"
$ws.Range("E39").Value = "11.72%"
$ws.Range("G39").Value = "12"
$ws.Range("D40").Value = "0.04736"
$ws.Range("E40").Value = "6.69%"
$ws.Range("G40").Value = "12"
$ws.Range("D41").Value = "0.007699"
$ws.Range("E41").Value = "4.40%"
$ws.Range("G41").Value = "12"
$ws.Range("E42").Value = "6.14%"
$ws.Range("G42").Value = "12"
$ws.Range("D43").Value = "0.007024"
$ws.Range("E43").Value = "-26.10%"
$ws.Range("G43").Value = "12"
$ws.Range("D44").Value = "0.002205"
$ws.Range("E44").Value = "-1.52%"
$ws.Range("G44").Value = "12"
$ws.Range("D45").Value = "0.009229"
$ws.Range("E45").Value = "14.84%"
$ws.Range("G45").Value = "12"
$ws.Range("D46").Value = "0.00006000"
$ws.Range("E46").Value = "-1.82%"
$ws.Range("G46").Value = "12"
$ws.Range("D47").Value = "0.00000000748"
$ws.Range("E47").Value = "-0.22%"
$ws.Range("G47").Value = "12"
$ws.Range("D48").Value = "5.826"
$ws.Range("E48").Value = "122.36%"
$ws.Range("G48").Value = "12"
$ws.Range("D49").Value = "0.002683"
$ws.Range("E49").Value = "34.06%"
$ws.Range("G49").Value = "12"
$ws.Range("D50").Value = "0.00002095"
$ws.Range("E50").Value = "-0.22%"
$ws.Range("G50").Value = "12"
$ws.Range("D51").Value = "0.0001995"
$ws.Range("E51").Value = "-0.22%"
$ws.Range("G51").Value = "12"
